$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ara" (Arabic) blacklisted-word rows entirely; this shifts
# the trailing blank rows up so the sheet ends at row 10 instead of 13.
$ws.Rows("8:10").Delete()

# Rename the "fra" language code to "spa" for the remaining French-text rows.
$ws.Range("A6").Value = "spa"
$ws.Range("A7").Value = "spa"

# The regenerated master data no longer carries the extra font/alignment
# styling on the data rows - reset A2:C7 back to the workbook's default style.
$ws.Range("A2:C7").Style = "Normal"

# Leave the selection where the author last left it when saving.
[void]$ws.Range("A8").Select()
